# Apply "More updates to the taxonomy" changes.
#
# Shape index map (slide 1, in Shapes collection order):
#  26 -> id=35 "TextBox 34"              ("Monitored database with / consensus")
#  29 -> id=54 "TextBox 53"              ("Replicated monitored database with consensus")
#  31 -> id=58 "Straight Arrow Connector 57" (glued connector feeding shape 35)
#  32 -> id=59 "TextBox 58"              ("Current state", rotated label)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) "Monitored database with" -> "Database with" (keep the trailing line break +
#    "consensus" run untouched) and shrink the box to its new autofit height.
$shp35 = $s.Shapes.Item(26)
$shp35.TextFrame.TextRange.Characters(1, 23).Text = "Database with"
$shp35.Height = 50.89220622440945

# 2) "Replicated monitored database with consensus" -> "Replicated database with consensus"
#    and shrink the box to its new autofit height.
$shp54 = $s.Shapes.Item(29)
$shp54.TextFrame.TextRange.Text = "Replicated database with consensus"
$shp54.Height = 72.70315170629921

# 3) Connector glued to shape 35 shortens to follow the smaller text box.
$shp58 = $s.Shapes.Item(31)
$shp58.Height = 35.020866141732284

# 4) "Current state" label rotates/moves slightly to stay aligned with the connector.
$shp59 = $s.Shapes.Item(32)
$shp59.Left = 462.2440157480315
$shp59.Top = 344.05070866141733
$shp59.Rotation = 19.757666666666665
